$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the text inside a Range with new text while forcing the
# native engine to keep it as its OWN run (instead of silently re-merging it
# into a neighbouring run that happens to share identical formatting).
#
# The engine normalises/merges adjacent runs that end up with identical
# run-properties at the moment a Range.Text assignment is made. Temporarily
# flipping a formatting flag (Bold) makes the freshly written text diverge
# from its neighbours, so the engine keeps it in its own <w:r>. Restoring the
# flag immediately afterwards (on a freshly re-fetched Range covering exactly
# the new text) removes the visual difference again without re-merging the
# run back into its neighbours.
# ---------------------------------------------------------------------------
function Set-RangeTextAsOwnRun($rng, [string]$newText) {
    $startPos = $rng.Start
    $rng.Font.Bold = 1
    $rng.Text = $newText
    $fixed = $d.Range($startPos, $startPos + $newText.Length)
    $fixed.Font.Bold = 0

    # The native engine coalesces any run immediately following the one we
    # just wrote if it happens to end up with identical run-properties (this
    # can silently swallow a short trailing run, e.g. a lone "." sentence
    # terminator, into the run we just split off). Touching the *entire*
    # following run's formatting (set + restore, over its full span so it
    # isn't itself split into pieces) makes the engine re-materialise it as
    # its own run again without altering its appearance.
    $afterPos = $startPos + $newText.Length
    $docEnd = $d.Content.End
    if ($afterPos -lt $docEnd) {
        $guard = $d.Range($afterPos, $afterPos + 1)
        if ($guard.Text -ne "") {
            $guardRun = $guard.Words(1)
            $guardRun.Font.Bold = 1
            $guardRun.Font.Bold = 0
        }
    }
}

# ---------------------------------------------------------------------------
# 1) Four occurrences of the "{{SEXO_7}}" merge field are renumbered to
#    "{{SEXO_19}}" by splicing "19" in place of the lone "7" digit.
# ---------------------------------------------------------------------------
$targets = @(
    "SE COMPROMETE A VENDER {{SEXO_7}} PROMITENTE {{SEXO_4}}",
    "EXPEDIR {{SEXO_7}} PROMITENTE {{SEXO_4}}",
    "REALIZAR EL REEMBOLSO {{SEXO_7}} PROMITENTE {{SEXO_4}}",
    "LIBERAR DE CARGA {{SEXO_7}} PROMITENTE {{SEXO_4}}"
)

foreach ($needle in $targets) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($needle)
    $digitIdx = $idx + $needle.IndexOf("{{SEXO_7}}") + ("{{SEXO_").Length
    $rng = $d.Range($digitIdx, $digitIdx + 1)
    Set-RangeTextAsOwnRun $rng "19"
}

# ---------------------------------------------------------------------------
# 2) The literal word " ÉSTE" following "... POR CUESTIONES AJENAS A" is
#    replaced by a new merge field. The clause tied to the PROMITENTE
#    {{SEXO_4}} paragraph becomes {{SEXO_17}}; the one tied to the
#    PROMITENTE {{SEXO_2}} paragraph becomes {{SEXO_16}}.
# ---------------------------------------------------------------------------
$esteReplacements = @(
    @{ Anchor = "LIBERAR DE CARGA {{SEXO_19}} PROMITENTE {{SEXO_4}}"; Field = "{{SEXO_17}}" },
    @{ Anchor = "LIBERAR DE CARGA {{SEXO_7}} PROMITENTE {{SEXO_2}}"; Field = "{{SEXO_16}}" }
)

foreach ($item in $esteReplacements) {
    $full = $d.Content.Text
    $anchorIdx = $full.IndexOf($item.Anchor)
    $searchFrom = $anchorIdx
    $needle = " POR CUESTIONES AJENAS A ÉSTE"
    $idx = $full.IndexOf($needle, $searchFrom)
    $steIdx = $idx + $needle.IndexOf("É")
    $rng = $d.Range($steIdx, $steIdx + 4)
    Set-RangeTextAsOwnRun $rng $item.Field
}
